$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 521, shifting the existing data (rows 521:540) down to (522:541)
$ws.Rows.Item(521).Insert()

# Populate the newly inserted row 521 with the new weekly price record
$ws.Range("A521").Value = 3
$ws.Range("B521").Value = "Femacal de La Calera"
$ws.Range("C521").Value = "Coquimbo"
$ws.Range("D521").NumberFormat = $ws.Range("D522").NumberFormat
$ws.Range("D521").Value = 45075
$ws.Range("E521").Value = 5
$ws.Range("F521").Value = 100112012
$ws.Range("G521").Value = "Espinaca"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 170
$ws.Range("K521").Value = 5500
$ws.Range("L521").Value = 6000
$ws.Range("M521").Value = 5765
$ws.Range("N521").Value = '$/docena de atados (3 kilos)'
$ws.Range("O521").Value = "Provincia de Quillota"
$ws.Range("P521").Value = 1922
$ws.Range("Q521").Value = 3
$ws.Range("R521").Value = "Hortaliza"
